$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 28, shifting existing rows 28:102 down to 29:103.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new weekly data point.
$ws.Range("A28").Value2 = 7
$ws.Range("B28").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C28").Value2 = "Ñuble"
$ws.Range("D28").Value2 = 44998
$ws.Range("E28").Value2 = 16
$ws.Range("F28").Value2 = "Fruta"
$ws.Range("G28").Value2 = 100108
$ws.Range("H28").Value2 = "Tropicales y subtropicales"
$ws.Range("I28").Value2 = 100108002
$ws.Range("J28").Value2 = "Mango"
$ws.Range("K28").Value2 = "Sin especificar"
$ws.Range("L28").Value2 = "Primera"
$ws.Range("M28").Value2 = 40
$ws.Range("N28").Value2 = 8000
$ws.Range("O28").Value2 = 8000
$ws.Range("P28").Value2 = 8000
$ws.Range("Q28").Value2 = "$/bandeja 4 kilos"
$ws.Range("R28").Value2 = "Perú"
$ws.Range("S28").Value2 = 2000
$ws.Range("T28").Value2 = 4
